$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.669.57"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "3.454.03"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'581.04"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").Value = "'145.78"
$ws.Range("E6").Value = "  +6.37%  "
$ws.Range("D7").Value = "3.455.56"
$ws.Range("E7").Value = "  +2.08%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.476"
$ws.Range("E9").Value = "  +1.53%  "
$ws.Range("D10").Value = "'7.62"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("E11").Value = "  +2.75%  "
$ws.Range("D12").Value = "'0.389"
$ws.Range("E12").Value = "  +2.18%  "
$ws.Range("D13").Value = "4.048.54"
$ws.Range("E13").Value = "  +2.18%  "
$ws.Range("E14").Value = "  +8.68%  "
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").Value = "3.454.45"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("D18").Value = "61.812.09"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").Value = "'6.25"
$ws.Range("E19").Value = "  +8.45%  "
$ws.Range("D20").Value = "'14.32"
$ws.Range("E20").Value = "  +3.64%  "
$ws.Range("D21").Value = "'9.57"
$ws.Range("E21").Value = "  +1.84%  "
$ws.Range("D22").Value = "'390.67"
$ws.Range("E22").Value = "  +3.68%  "
$ws.Range("D23").Value = "'0.566"
$ws.Range("E23").Value = "  +2.95%  "
$ws.Range("D24").Value = "'73.83"
$ws.Range("E24").Value = "  +3.88%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "'5.78"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("D28").Value = "3.595.88"
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("D29").Value = "'0.182"
$ws.Range("E29").Value = "  +1.59%  "
$ws.Range("D30").Value = "'7.62"
$ws.Range("E30").Value = "  +2.61%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").Value = "'8.18"
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'2.19"
$ws.Range("E33").Value = "  +2.14%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.46"
$ws.Range("E34").Value = "  -11.13%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "'24.08"
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("D37").Value = "3.484.98"
$ws.Range("E37").Value = "  +2.22%  "
$ws.Range("D38").Value = "'7.02"
$ws.Range("E38").Value = "  +2.90%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "'5.14"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "'1.56"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").Value = "'166.61"
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("D42").Value = "'0.0783"
$ws.Range("E42").Value = "  +2.69%  "
$ws.Range("D43").Value = "'27.72"
$ws.Range("E43").Value = "  +8.69%  "
$ws.Range("D44").Value = "'0.807"
$ws.Range("E44").Value = "  +3.93%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'42.49"
$ws.Range("E45").Value = "  +1.87%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'4.50"
$ws.Range("E46").Value = "  +3.72%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("D50").Value = "2.572.56"
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("D51").Value = "'6.94"
$ws.Range("E51").Value = "  +2.22%  "
